$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the indicator_last_checked timestamp for row 2 (AAGIY) to reflect
# the latest fetch run.
$ws.Range("C2").Value = "2025-09-24T19:34:13.184636"

# The mock indicator data that had been fetched for rows 3-11 (AAPL, ABBV,
# ADBE, ADI, ADSK, ADYEY, AEM, AJG, ALL) turned out to be bogus (network/DNS
# resolution was unavailable, so a mock-data fallback had populated these
# rows). Blank out columns B..V for those rows, leaving only the Ticker in
# column A - matching the blank placeholder rows already present further
# down the sheet (row 12 onward) that are still awaiting real data.
#
# A leading-apostrophe assignment clears each cell to an empty (blank) text
# value rather than deleting it outright, matching the existing blank rows;
# re-applying the "Normal" style afterwards strips the text-quote-prefix
# formatting that the apostrophe entry would otherwise leave behind, so the
# cells end up plain/unstyled just like the other blank rows.
$ws.Range("B3:V11").Value = "'"
$ws.Range("B3:V11").Style = "Normal"
